$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# Remove the two "separator" rows that only contained a note in column J
# (original rows 24 and 26). Deleting row 24 shifts everything up, which
# makes the former row 26 become row 25 - delete that next.
# ------------------------------------------------------------------
$ws.Rows(24).Delete()
$ws.Rows(25).Delete()

# ------------------------------------------------------------------
# Mark the "Parceiro Exsat" rows (now rows 20-25) as implemented by
# applying the same "Bom" (Good) cell style used elsewhere in the sheet
# and adding the "Implementado" note in column K.
# ------------------------------------------------------------------
$ws.Range("A20:J25").Style = "Bom"
$ws.Range("C20:D25").NumberFormat = "0%"
$ws.Range("E20:E25").NumberFormat = "0.00%"
$ws.Range("G20:H25").NumberFormat = "0.00%"

$ws.Range("K20").Value = "Implementado"
$ws.Range("K21").Value = "Implementado"
$ws.Range("K22").Value = "Implementado"
$ws.Range("K23").Value = "Implementado"
$ws.Range("K24").Value = "Implementado"
$ws.Range("K25").Value = "Implementado"
$ws.Range("K20:K25").Style = "Bom"

# ------------------------------------------------------------------
# Restore the view: scroll the frozen pane back to the top and move the
# active selection to I25.
# ------------------------------------------------------------------
$aw = $excel.ActiveWindow
$aw.ScrollRow = 2
$aw.ScrollColumn = 1
$ws.Range("I25").Select()
